$p = $ppt.ActivePresentation

# Add a new "Title and Content" slide (CustomLayout index 2, i.e.
# ppt/slideLayouts/slideLayout2.xml) at the very end of the deck.
$count = $p.Slides.Count
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s = $p.Slides.AddSlide($count + 1, $layout)

# Title: "GitKraken" + " Graphical View" (two runs, like the source deck),
# centered.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "GitKraken"
[void]$title.InsertAfter(" Graphical View")
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.ParagraphFormat.Alignment = 2  # ppAlignCenter

# Give the body placeholder its usual authoring name; leave its text empty.
$s.Shapes.Item(2).Name = "Content Placeholder 6"
